$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 370) holds a date value (serial 45190 = 2023-09-21)
# that was bumped forward by two days (to serial 45192 = 2023-09-23) for every row.
$lastRow = 370
$rng = $ws.Range("C2:C" + $lastRow)
$rng.Value = 45192
